# Ajout croquis - guide organisateur
# Adds the "Senneterre" dinner/souper option (mirroring "Malartic"), replaces the
# breakfast description text, replaces the DINER "Sam" row with the lunch-box
# text (already used elsewhere), and updates the sheet/window selection state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# DEJEUNER sheet: new breakfast wording (row 2)
# ---------------------------------------------------------------------------
$wsDej = $wb.Worksheets.Item("DEJEUNER")
$wsDej.Range("B2").Value = "Œufs, bacon, jambon, saucisse, pommes de terre rôties, crêpe avec sirop, céréales froides, gruau, banane, orange ou pomme, breuvage (jus,, lait, café)"
$wsDej.Range("C2").Value = "Eggs, bacon, ham, sausage,roasted potatoes, pancake with syrup, cold cereal, oatmeal, banana, orange or apple, beverage  juice, milk, coffee"
$wsDej.Rows.Item(2).RowHeight = 17

# ---------------------------------------------------------------------------
# DINER sheet: "Sam" row (row 6) now carries the lunch-box description
# ---------------------------------------------------------------------------
$wsDin = $wb.Worksheets.Item("DINER")
$wsDin.Range("B6").Value = "**Boîte à lunch**<br/> Wrap à la viande froide, crudités, fromage, salade de macaroni, biscuits, orange, jus"
$wsDin.Range("C6").Value = "**Lunch box**<br/> Cold cuts wrap, raw vegetables, cheese, macaroni salad, cookies, orange, juice"
$wsDin.Range("B6").HorizontalAlignment = -4131
$wsDin.Range("B6").WrapText = $true
$wsDin.Range("C6").HorizontalAlignment = -4131
$wsDin.Range("C6").WrapText = $true
$wsDin.Rows.Item(6).RowHeight = 17

# ---------------------------------------------------------------------------
# SOUPER sheet: "Ven" row (row 7) becomes the new "Repas à Senneterre" entry
# ---------------------------------------------------------------------------
$wsSou = $wb.Worksheets.Item("SOUPER")
$wsSou.Range("B7").Value = "Repas à Senneterre"
$wsSou.Range("C7").Value = "Dinner at Senneterre"
$wsSou.Rows.Item(7).RowHeight = 17

# ---------------------------------------------------------------------------
# Window / selection state: DEJEUNER becomes the active tab, each sheet keeps
# track of its own last-used selection.
# ---------------------------------------------------------------------------
$wsDin.Activate()
$wsDin.Range("B14").Select()

$wsSou.Activate()
$wsSou.Range("B8").Select()

$wsDej.Activate()
$wsDej.Range("C3").Select()
